# Weekly refresh of the "Perejil" price sheet.
#
# The new reporting week's data lands at the top of the table (rows 26-27,
# the 13th date-pair for this market/category) and every older date-pair
# that used to occupy rows 26-99 is pushed down by one pair (2 rows), with
# the oldest pair (previously rows 98-99) spilling into two brand-new rows
# (100-101) at the bottom. Columns A-R move together as a unit for every
# pair; only the newest pair's date (D26/D27) is genuinely new data - it is
# one week later than what is now D28/D29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 26
$lastDataRow  = 99
$lastCol      = 18   # A..R
$dateCol      = 4    # D
$newDate      = 44453

# Walk bottom-up so we never clobber a source row before it has been read.
for ($r = $lastDataRow; $r -ge $firstDataRow; $r--) {
    $destRow = $r + 2
    for ($c = 1; $c -le $lastCol; $c++) {
        $srcCell = $ws.Cells.Item($r, $c)
        $dstCell = $ws.Cells.Item($destRow, $c)
        $dstCell.Value2 = $srcCell.Value2
    }
    # Column D carries an explicit date number format; carry it along so the
    # newly-created rows 100/101 render as dates too.
    $ws.Cells.Item($destRow, $dateCol).NumberFormat = $ws.Cells.Item($r, $dateCol).NumberFormat
}

# The top two rows keep all of their original data except the date, which
# advances to the new reporting week.
$ws.Cells.Item($firstDataRow, $dateCol).Value2 = $newDate
$ws.Cells.Item($firstDataRow + 1, $dateCol).Value2 = $newDate
